$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.245.15'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '2.168.18'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.12'
$ws.Range("E5").Value = '  +5.97%  '
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.87'
$ws.Range("E7").Value = '  +1.32%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.578'
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.63'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '2.497.07'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("E15").Value = '  -2.23%  '
$ws.Range("D16").Value = '2.172.03'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.765'
$ws.Range("E17").Value = '  -2.69%  '
$ws.Range("D18").Value = '42.139.04'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.43'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.19'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.14'
$ws.Range("E23").Value = '  +3.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.41'
$ws.Range("E24").Value = '  -6.04%  '
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.38'
$ws.Range("E26").Value = '  -3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.36'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.14'
$ws.Range("E29").Value = '  -1.97%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.37'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.38'
$ws.Range("E31").Value = '  +9.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.93'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0808'
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("E34").Value = '  -3.89%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("E36").Value = '  +1.89%  '
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("E39").Value = '  -2.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.72'
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '59.10'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.12'
$ws.Range("E43").Value = '  -4.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.70'
$ws.Range("E44").Value = '  +4.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.466'
$ws.Range("E45").Value = '  +11.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0969'
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.22'
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.41'
$ws.Range("E48").Value = '  +8.57%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("E51").Value = '  +0.32%  '
